# Generate Report for Handoff
# Updates the "Latest Handoff Date"/"Latest Handoff Datetime" timestamps to reflect
# the newly generated handoff report for e7e9493d-8e09-4804-8747-d514cf7f7ff3.

$wb = $excel.ActiveWorkbook

# Overview sheet: Latest Handoff Date for e7e9493d-8e09-4804-8747-d514cf7f7ff3.md (row 7)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D7").Value = "2016-03-25 09:09:26"

# zh-cn sheet: Latest Handoff Datetime for e7e9493d-8e09-4804-8747-d514cf7f7ff3 (row 7)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E7").Value = "2016-03-25 09:09:21"

# de-de sheet: Latest Handoff Datetime for e7e9493d-8e09-4804-8747-d514cf7f7ff3 (row 7)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E7").Value = "2016-03-25 09:09:26"
